$wb = $excel.ActiveWorkbook
try {
  Write-Host "Path:" $wb.Path
  Write-Host "FullName:" $wb.FullName
} catch {
  Write-Host "err" $_
}
